{"js": "// update semaine 7 v7.0\n// Append a new bullet item (\"EA -> fixer objectifs pour ann\u00e9e prochaine.\")\n// at the very end of the document, continuing the same bulleted list\n// (style \"Paragraphedeliste\" / numId 1) as the paragraph before it.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\n\n// Inserting right after the last (list) paragraph continues its list/style,\n// exactly like pressing Enter at the end of that bullet in Word.\nlastParagraph.insertParagraph(\n  \"EA -> fixer objectifs pour ann\u00e9e prochaine.\",\n  Word.InsertLocation.after\n);\n\nawait context.sync();\n", "ps1": "# update semaine 7 v7.0\n# Append a new bullet item (\"EA -> fixer objectifs pour ann\u00e9e prochaine.\")\n# to the end of the document, continuing the same bulleted list\n# (Paragraphedeliste / numId 1) as the preceding item.\n\n$d = $word.ActiveDocument\n\n# Last paragraph in the document body (\"... lors d'un packaging...\")\n$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)\n\n# Insert a new paragraph right after it; Word continues the enclosing\n# list/style automatically (same as pressing Enter at the end of the item).\n$lastPara.Range.InsertParagraphAfter()\n\n$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)\n$newPara.Range.Text = \"EA -> fixer objectifs pour ann\u00e9e prochaine.\"\n"}
